{"js": "const body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Change 1: \"Depois a frota passou em Sic\u00edlia...\" ->\n//           \"Depois de Caffa, a frota passou em Sic\u00edlia...\"\n// ---------------------------------------------------------------------\nconst seg1 = body.search(\n  \"Depois a frota passou em Sic\u00edlia que ficou por tr\u00eas semanas\",\n  { matchCase: true }\n);\nseg1.load(\"text\");\nawait context.sync();\n\nif (seg1.items.length > 0) {\n  seg1.items[0].insertText(\n    \"Depois de Caffa, a frota passou em Sic\u00edlia que ficou por tr\u00eas semanas\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// The existing \"_GoBack\" bookmark (originally right after the first\n// \"Marselha\") moves to sit right after the newly inserted \"Caffa,\".\n// Remove it from its old spot first so the name stays unique.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-insert \"_GoBack\" as a collapsed bookmark right before \" a frota\n// passou em Sic\u00edlia\" (i.e. immediately after \"Caffa,\").\nconst insertionAnchor = body.search(\" a frota passou em Sic\u00edlia que ficou por tr\u00eas semanas\", {\n  matchCase: true,\n});\ninsertionAnchor.load(\"text\");\nawait context.sync();\n\nif (insertionAnchor.items.length > 0) {\n  const collapsedStart = insertionAnchor.items[0].getRange(\"Start\");\n  collapsedStart.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Change 2: \"tanto que a frota foi expulsar de l\u00e1, e s\u00f3 conseguir \" ->\n//           \"at\u00e9 a frota for expulsa de l\u00e1,, e s\u00f3 conseguir \"\n// ---------------------------------------------------------------------\nconst seg2 = body.search(\n  \"tanto que a frota foi expulsar de l\u00e1, e s\u00f3 conseguir \",\n  { matchCase: true }\n);\nseg2.load(\"text\");\nawait context.sync();\n\nif (seg2.items.length > 0) {\n  seg2.items[0].insertText(\n    \"at\u00e9 a frota for expulsa de l\u00e1,, e s\u00f3 conseguir \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Standard Word COM constants (defined explicitly; this host does not\n# pre-seed the usual $wdXxx globals).\n$wdReplaceOne    = 1\n$wdFindContinue  = 1\n$wdCollapseEnd   = 0\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Change 1: \"Depois a frota passou em Sic\u00edlia...\" ->\n#           \"Depois de Caffa, a frota passou em Sic\u00edlia...\"\n# ---------------------------------------------------------------------\n$rng1 = $d.Content\n$rng1.Find.Execute(\n    \"Depois a frota passou em Sic\u00edlia que ficou por tr\u00eas semanas\",\n    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \"Depois de Caffa, a frota passou em Sic\u00edlia que ficou por tr\u00eas semanas\",\n    $wdReplaceOne\n) | Out-Null\n\n# The \"_GoBack\" bookmark that used to sit right after the first\n# \"Marselha\" moves to sit right after the newly-inserted \"Caffa,\".\n# Remove it from its old location first so the name stays unique.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Locate the just-inserted \"Caffa,\" (unique because it is followed by\n# \" a frota\", unlike the pre-existing \"Caffa,\" earlier in the document)\n# and drop a collapsed bookmark right after it.\n$rngCaffa = $d.Content\n$rngCaffa.Find.Execute(\"Caffa, a frota\") | Out-Null\n$rngCaffa.SetRange($rngCaffa.Start, $rngCaffa.Start + 6)\n$rngCaffa.Collapse($wdCollapseEnd)\n$d.Bookmarks.Add(\"_GoBack\", $rngCaffa) | Out-Null\n\n# ---------------------------------------------------------------------\n# Change 2: \"tanto que a frota foi expulsar de l\u00e1, e s\u00f3 conseguir \" ->\n#           \"at\u00e9 a frota for expulsa de l\u00e1,, e s\u00f3 conseguir \"\n# ---------------------------------------------------------------------\n$rng2 = $d.Content\n$rng2.Find.Execute(\n    \"tanto que a frota foi expulsar de l\u00e1, e s\u00f3 conseguir \",\n    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \"at\u00e9 a frota for expulsa de l\u00e1,, e s\u00f3 conseguir \",\n    $wdReplaceOne\n) | Out-Null\n"}
